# Update Mainboard pick-and-place table to the new board revision:
#  - the CFG connector was removed from this revision
#  - remaining parts were re-sequenced (LED/X/Y/Z/DISP/DB1 first, then the
#    XY-XH2.54 headers C1-C4/TOP)
#  - the LED footprint moved slightly on the new board

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (CFG) - board revision no longer includes it
$ws.Rows.Item(13).Delete()

# Rewrite rows 2-12 with the new board revision data (reordered + updated LED position)
# Row 2: LED
$ws.Range("A2").Value = "LED"
$ws.Range("B2").Value = "ZX-XH2.54-3PZZ"
$ws.Range("C2").Value = "CONN-TH_3P-P2.54_2501S-3P"
$ws.Range("D2").Value = "80.899mm"
$ws.Range("E2").Value = "-54.102mm"
$ws.Range("F2").Value = "80.899mm"
$ws.Range("G2").Value = "-54.102mm"
$ws.Range("H2").Value = "78.359mm"
$ws.Range("I2").Value = "-54.102mm"
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = "T"
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = "No"
$ws.Range("N2").Value = "ZX-XH2.54-3PZZ"

# Row 3: X
$ws.Range("A3").Value = "X"
$ws.Range("B3").Value = "ZX-XH2.54-3PZZ"
$ws.Range("C3").Value = "CONN-TH_3P-P2.54_2501S-3P"
$ws.Range("D3").Value = "80.899mm"
$ws.Range("E3").Value = "-43.942mm"
$ws.Range("F3").Value = "80.899mm"
$ws.Range("G3").Value = "-43.942mm"
$ws.Range("H3").Value = "83.439mm"
$ws.Range("I3").Value = "-43.942mm"
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = "T"
$ws.Range("L3").Value = 180
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "ZX-XH2.54-3PZZ"

# Row 4: Y
$ws.Range("A4").Value = "Y"
$ws.Range("B4").Value = "ZX-XH2.54-3PZZ"
$ws.Range("C4").Value = "CONN-TH_3P-P2.54_2501S-3P"
$ws.Range("D4").Value = "80.899mm"
$ws.Range("E4").Value = "-37.465mm"
$ws.Range("F4").Value = "80.899mm"
$ws.Range("G4").Value = "-37.465mm"
$ws.Range("H4").Value = "83.439mm"
$ws.Range("I4").Value = "-37.465mm"
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = "T"
$ws.Range("L4").Value = 180
$ws.Range("M4").Value = "No"
$ws.Range("N4").Value = "ZX-XH2.54-3PZZ"

# Row 5: Z
$ws.Range("A5").Value = "Z"
$ws.Range("B5").Value = "ZX-XH2.54-3PZZ"
$ws.Range("C5").Value = "CONN-TH_3P-P2.54_2501S-3P"
$ws.Range("D5").Value = "80.899mm"
$ws.Range("E5").Value = "-30.988mm"
$ws.Range("F5").Value = "80.899mm"
$ws.Range("G5").Value = "-30.988mm"
$ws.Range("H5").Value = "83.439mm"
$ws.Range("I5").Value = "-30.988mm"
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = "T"
$ws.Range("L5").Value = 180
$ws.Range("M5").Value = "No"
$ws.Range("N5").Value = "ZX-XH2.54-3PZZ"

# Row 6: DISP
$ws.Range("A6").Value = "DISP"
$ws.Range("B6").Value = "ZX-XH2.54-4PZZ"
$ws.Range("C6").Value = "CONN-TH_4P-P2.50_4PIN"
$ws.Range("D6").Value = "43.18mm"
$ws.Range("E6").Value = "-19.304mm"
$ws.Range("F6").Value = "43.18mm"
$ws.Range("G6").Value = "-19.304mm"
$ws.Range("H6").Value = "43.18mm"
$ws.Range("I6").Value = "-15.554mm"
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = "T"
$ws.Range("L6").Value = 270
$ws.Range("M6").Value = "No"
$ws.Range("N6").Value = "ZX-XH2.54-4PZZ"

# Row 7: DB1
$ws.Range("A7").Value = "DB1"
$ws.Range("B7").Value = "Pico Pi RP2040"
$ws.Range("C7").Value = "YD-RP2040"
$ws.Range("D7").Value = "61.468mm"
$ws.Range("E7").Value = "-37.338mm"
$ws.Range("F7").Value = "52.578mm"
$ws.Range("G7").Value = "-13.208mm"
$ws.Range("H7").Value = "52.578mm"
$ws.Range("I7").Value = "-13.208mm"
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = "T"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "No"
$ws.Range("N7").Value = "YD-RP2040"

# Row 8: C1
$ws.Range("A8").Value = "C1"
$ws.Range("B8").Value = "XY-XH2.54-2A11"
$ws.Range("C8").Value = "HDR-TH_2P-P2.50-V-F-1"
$ws.Range("D8").Value = "43.307mm"
$ws.Range("E8").Value = "-31.369mm"
$ws.Range("F8").Value = "43.307mm"
$ws.Range("G8").Value = "-31.369mm"
$ws.Range("H8").Value = "44.577mm"
$ws.Range("I8").Value = "-31.369mm"
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = "T"
$ws.Range("L8").Value = 180
$ws.Range("M8").Value = "No"
$ws.Range("N8").Value = "XY-XH2.54-2A11"

# Row 9: C2
$ws.Range("A9").Value = "C2"
$ws.Range("B9").Value = "XY-XH2.54-2A11"
$ws.Range("C9").Value = "HDR-TH_2P-P2.50-V-F-1"
$ws.Range("D9").Value = "43.434mm"
$ws.Range("E9").Value = "-38.989mm"
$ws.Range("F9").Value = "43.434mm"
$ws.Range("G9").Value = "-38.989mm"
$ws.Range("H9").Value = "44.704mm"
$ws.Range("I9").Value = "-38.989mm"
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = "T"
$ws.Range("L9").Value = 180
$ws.Range("M9").Value = "No"
$ws.Range("N9").Value = "XY-XH2.54-2A11"

# Row 10: C3
$ws.Range("A10").Value = "C3"
$ws.Range("B10").Value = "XY-XH2.54-2A11"
$ws.Range("C10").Value = "HDR-TH_2P-P2.50-V-F-1"
$ws.Range("D10").Value = "43.434mm"
$ws.Range("E10").Value = "-46.863mm"
$ws.Range("F10").Value = "43.434mm"
$ws.Range("G10").Value = "-46.863mm"
$ws.Range("H10").Value = "44.704mm"
$ws.Range("I10").Value = "-46.863mm"
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = "T"
$ws.Range("L10").Value = 180
$ws.Range("M10").Value = "No"
$ws.Range("N10").Value = "XY-XH2.54-2A11"

# Row 11: C4
$ws.Range("A11").Value = "C4"
$ws.Range("B11").Value = "XY-XH2.54-2A11"
$ws.Range("C11").Value = "HDR-TH_2P-P2.50-V-F-1"
$ws.Range("D11").Value = "43.434mm"
$ws.Range("E11").Value = "-54.737mm"
$ws.Range("F11").Value = "43.434mm"
$ws.Range("G11").Value = "-54.737mm"
$ws.Range("H11").Value = "44.704mm"
$ws.Range("I11").Value = "-54.737mm"
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = "T"
$ws.Range("L11").Value = 180
$ws.Range("M11").Value = "No"
$ws.Range("N11").Value = "XY-XH2.54-2A11"

# Row 12: TOP
$ws.Range("A12").Value = "TOP"
$ws.Range("B12").Value = "XY-XH2.54-2A11"
$ws.Range("C12").Value = "HDR-TH_2P-P2.50-V-F-1"
$ws.Range("D12").Value = "43.307mm"
$ws.Range("E12").Value = "-62.611mm"
$ws.Range("F12").Value = "43.307mm"
$ws.Range("G12").Value = "-62.611mm"
$ws.Range("H12").Value = "44.577mm"
$ws.Range("I12").Value = "-62.611mm"
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = "T"
$ws.Range("L12").Value = 180
$ws.Range("M12").Value = "No"
$ws.Range("N12").Value = "XY-XH2.54-2A11"
